# This script reproduces the edit captured in the diff:
#   - Column C's header on the "fromagerie" sheet is renamed from
#     "Photo" to "Image_Path" (C1). That is the only cell whose text
#     actually changes anywhere in the workbook.
#   - Because "Photo" then becomes an unused shared string, Excel drops
#     it from sharedStrings.xml on save, which shifts every subsequent
#     shared-string index down by one. That index shift is exactly what
#     produces all of the <v>N</v> -> <v>N-1</v> churn seen throughout
#     the raw OOXML diff (other sheets' C1 already said "Image_Path",
#     just via a different, soon-to-be-renumbered shared string index)
#     -- no other cell's actual content changes.
#   - The selection/active cell on a couple of sheets simply reflects
#     where the user's cursor ended up after making the edit.

$wb = $excel.ActiveWorkbook

# Update header cell C1 on "fromagerie" from "Photo" to "Image_Path"
$wsFromagerie = $wb.Worksheets.Item("fromagerie")
$wsFromagerie.Range("C1").Value = "Image_Path"

# "apiculture" already reads "Image_Path" in C1; re-assert it so the
# workbook is fully self-consistent even if that ever changes upstream.
$wsApiculture = $wb.Worksheets.Item("apiculture")
$wsApiculture.Range("C1").Value = "Image_Path"

# Reflect the resulting cursor/selection positions recorded in the diff
$wsFromagerie.Activate()
$wsFromagerie.Range("B12").Select()

$wsApiculture.Activate()
$wsApiculture.Range("C1").Select()

# Restore original active sheet (first sheet was tabSelected in the file)
$wsFromagerie.Activate()
